# Atualizando o andamento das SARs.

$wb = $excel.ActiveWorkbook

# --- Sheet "Histórico de SARs": update status + last-change date for all 3 SARs ---
$ws1 = $wb.Worksheets.Item("Histórico de SARs")

$ws1.Range("C2").Value = "Aprovada para Resolução"
$ws1.Range("C3").Value = "Aprovada para Resolução"
$ws1.Range("C4").Value = "Aprovada para Resolução"

$ws1.Range("F2").Value = "5/21/2015"
$ws1.Range("F3").Value = "5/21/2015"
$ws1.Range("F4").Value = "5/21/2015"

# Widen column C so the longer status text fits.
$ws1.Columns.Item(3).ColumnWidth = 23.5

# --- Sheet "Legenda do Documento": merge "Efetivada"/"Encerrada" legend rows ---
$ws2 = $wb.Worksheets.Item("Legenda do Documento")

$ws2.Range("C5").Value = "Efetivada e Encerrada"
$ws2.Rows.Item(6).Delete()

$ws2.Range("C9").Select()

# Re-select sheet 1 last so it stays the active tab (tabSelected).
$ws1.Range("F4").Select()
